$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.381.23'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.846.54'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6321'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07597'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2974'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.42'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '2.280.10'
$ws.Range('E11').Value = '  +22.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07720'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.990'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6864'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '2.398.18'
$ws.Range('E15').Value = '  +12.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.92'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009877'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.160'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').Value = '29.392.62'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '231.70'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.50'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.592'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.64'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.465'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.474'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05800'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.253'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.018'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.862'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7175'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = '1.250.21'
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.795'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01802'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').Value = '2.343.32'
$ws.Range('E41').Value = '  +15.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9050'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.105'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9994'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.38'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.38'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.288'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.153'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4008'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.694'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05747'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.17%  '
